# Generate Report for Handback
# The localization handback for f674905b-9063-4abe-af0b-b694f7aa8537.md has
# completed: its status moves from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, the "Latest Handback
# DateTime" stamps are refreshed for both target locales, and the stale
# "version not latest" error is cleared now that the handback succeeded.

$wb = $excel.ActiveWorkbook

$status_done = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the f674905b-...md file ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = $status_done   # zh-cn status
$ov.Range("F3").Value = $status_done   # de-de status

# --- zh-cn sheet: row 3 is the f674905b-...md file ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $status_done              # Status
$zh.Range("K3").Value = "2016-08-21 08:53:57"     # Latest Handback DateTime
$zh.Range("P3").Value = ""                        # Error Detail cleared
$zh.Columns.Item(16).ColumnWidth = 12.8            # Error Detail col auto-shrinks

# --- de-de sheet: row 3 is the f674905b-...md file ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $status_done              # Status
$de.Range("K3").Value = "2016-08-21 08:54:07"     # Latest Handback DateTime
$de.Range("P3").Value = ""                        # Error Detail cleared
$de.Columns.Item(16).ColumnWidth = 12.8            # Error Detail col auto-shrinks
